$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.713.26"
$ws.Range("E2").Value = "  +2.15%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.892.17"
$ws.Range("E3").Value = "  +0.77%  "
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.16"
$ws.Range("E5").Value = "  +0.68%  "
$ws.Range("E6").Value = "  +0.15%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4926"
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2957"
$ws.Range("E8").Value = "  +0.86%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06799"
$ws.Range("E9").Value = "  +2.65%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.889.89"
$ws.Range("E10").Value = "  +0.64%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "17.23"
$ws.Range("E11").Value = "  +3.97%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07245"
$ws.Range("E12").Value = "  +0.52%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "90.92"
$ws.Range("E13").Value = "  +5.19%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6787"
$ws.Range("E14").Value = "  +1.73%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.045"
$ws.Range("E15").Value = "  +2.45%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "30.668.04"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000007988"
$ws.Range("E17").Value = "  +2.04%  "
$ws.Range("E18").Value = "  +0.10%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.16"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.131.19"
$ws.Range("E20").Value = "  +0.45%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.000"
$ws.Range("E21").Value = "  -0.07%  "
$ws.Range("E22").Value = "  +0.66%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "189.35"
$ws.Range("E23").Value = "  +32.78%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.148"
$ws.Range("E24").Value = "  +4.80%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.405"
$ws.Range("E25").Value = "  +2.75%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "155.86"
$ws.Range("E26").Value = "  +2.24%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.09"
$ws.Range("E27").Value = "  +12.51%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.900"
$ws.Range("E28").Value = "  -0.03%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.400"
$ws.Range("E29").Value = "  +1.17%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.345"
$ws.Range("E30").Value = "  +3.18%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09075"
$ws.Range("E31").Value = "  +3.16%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.013"
$ws.Range("E32").Value = "  +0.54%  "
$ws.Range("E33").Value = "  +2.53%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7513"
$ws.Range("E34").Value = "  +4.32%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.109"
$ws.Range("E35").Value = "  -0.40%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.775"
$ws.Range("E36").Value = "  +4.33%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01839"
$ws.Range("E37").Value = "  -0.59%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.688"
$ws.Range("E38").Value = "  +0.21%  "
$ws.Range("E39").Value = "  -0.84%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9375"
$ws.Range("E40").Value = "  +0.81%  "
$ws.Range("E41").Value = "  +4.48%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "105.37"
$ws.Range("E42").Value = "  +2.01%  "
$ws.Range("E43").Value = "  +0.23%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.769"
$ws.Range("E44").Value = "  -0.27%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.589"
$ws.Range("E45").Value = "  +2.75%  "
$ws.Range("E46").Value = "  +4.96%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05868"
$ws.Range("E47").Value = "  +2.96%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.735"
$ws.Range("E48").Value = "  +4.97%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.422"
$ws.Range("E49").Value = "  +5.98%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.3935"
$ws.Range("E50").Value = "  +3.98%  "
$ws.Range("E51").Value = "  +2.35%  "
